$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture comment texts that live on rows which are about to shift down ---
$commentH142 = $ws.Range('H142').Comment.Text()
$commentH147 = $ws.Range('H147').Comment.Text()

# Remove the old comments before the row insert shuffles everything
$ws.Range('H142').Comment.Delete()
$ws.Range('H147').Comment.Delete()

# --- Insert a new row at 142, pushing existing rows 142-164 down to 143-165 ---
$ws.Rows("142").Insert()

# --- Populate the newly inserted row 142 ---
# Write order matters for shared-string allocation order (B, H, then D)
$ws.Range("A142").Value = "AESutils"
$ws.Range("B142").Value = "getelemsGUI"
$ws.Range("H142").Value = "interactively choose elements for quant or plots"
$ws.Range("D142").Value = "AESquantparams"

# --- Re-create the shifted comments at their new locations ---
$ws.Range('H143').AddComment($commentH142)
$ws.Range('H148').AddComment($commentH147)

# --- Add the new comment describing the D142 cell ---
$ws.Range('D142').AddComment("Lines only available if specified in AESquantparams")

# --- Update the filter-database defined name to cover the new last row ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Functions!_FilterDatabase") {
        $n.RefersTo = "=Functions!`$A`$1:`$I`$165"
    }
}

# --- Restore selection to reflect where editing ended up ---
$ws.Range("D145").Select()
